$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header columns to English/snake_case field names
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# Title-case the "de"/"el" connector in municipality / state names
$ws.Range("B6").Value = "Amatenango De La Frontera"
$ws.Range("A15").Value = "Ciudad De México"
$ws.Range("A19").Value = "Estado De México"
$ws.Range("B19").Value = "Almoloya De Juárez"
$ws.Range("B28").Value = "Alcozauca De Guerrero"
$ws.Range("B29").Value = "Atoyac De Álvarez"
$ws.Range("B30").Value = "Chilapa De Álvarez"
$ws.Range("B31").Value = "Coyuca De Catalán"
$ws.Range("B37").Value = "Técpan De Galeana"
$ws.Range("B39").Value = "Atotonilco El Grande"
$ws.Range("B44").Value = "Tenango De Doria"
$ws.Range("B48").Value = "Ojuelos De Jalisco"
$ws.Range("B62").Value = "Acatlán De Pérez Figueroa"
$ws.Range("B64").Value = "Nejapa De Madero"
$ws.Range("B65").Value = "Oaxaca De Juárez"
$ws.Range("B80").Value = "San Salvador El Verde"
$ws.Range("B86").Value = "San Ciro De Acosta"
$ws.Range("B92").Value = "Hueyapan De Ocampo"

# Remove the trailing metadata/footer rows (104-108); row 103 is blank/unused
$ws.Rows("104:108").Delete()
